$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.209.32'
$ws.Range("E2").Value = '  +1.17%  '
$ws.Range("D3").Value = '1.835.15'
$ws.Range("E3").Value = '  +1.06%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.013'
$ws.Range("E4").Value = '  +1.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.92'
$ws.Range("E5").Value = '  +1.23%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.011'
$ws.Range("E6").Value = '  +1.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4711'
$ws.Range("E7").Value = '  +0.30%  '
$ws.Range("E8").Value = '  -0.28%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07423'
$ws.Range("E9").Value = '  +0.60%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8825'
$ws.Range("E10").Value = '  +1.35%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.45'
$ws.Range("E11").Value = '  +0.08%  '
$ws.Range("D12").Value = '1.830.16'
$ws.Range("E12").Value = '  +0.52%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07329'
$ws.Range("E13").Value = '  +3.47%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.478'
$ws.Range("E14").Value = '  +2.04%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.80'
$ws.Range("E15").Value = '  +0.65%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.570'
$ws.Range("E16").Value = '  +0.99%  '
$ws.Range("E17").Value = '  +1.14%  '
$ws.Range("E18").Value = '  +0.75%  '
$ws.Range("E19").Value = '  +0.90%  '
$ws.Range("E20").Value = '  +0.37%  '
$ws.Range("D21").Value = '27.230.15'
$ws.Range("E21").Value = '  +1.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.310'
$ws.Range("E22").Value = '  -0.63%  '
$ws.Range("E23").Value = '  +1.31%  '
$ws.Range("D24").Value = '2.054.40'
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.904'
$ws.Range("E25").Value = '  +0.59%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '152.43'
$ws.Range("E26").Value = '  +0.27%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.63'
$ws.Range("E27").Value = '  +1.34%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.166'
$ws.Range("E28").Value = '  -1.33%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.282'
$ws.Range("E29").Value = '  -0.47%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '117.84'
$ws.Range("E30").Value = '  +2.12%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08927'
$ws.Range("E31").Value = '  +0.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7600'
$ws.Range("E32").Value = '  -0.97%  '
$ws.Range("E33").Value = '  +0.83%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.550'
$ws.Range("E34").Value = '  +1.46%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.941'
$ws.Range("E35").Value = '  +0.72%  '
$ws.Range("E36").Value = '  +1.03%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.106'
$ws.Range("E37").Value = '  +0.65%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05341'
$ws.Range("E38").Value = '  +1.56%  '
$ws.Range("E39").Value = '  +0.16%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.011'
$ws.Range("E40").Value = '  +2.24%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.424'
$ws.Range("E41").Value = '  +2.93%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '7.344'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5356'
$ws.Range("E43").Value = '  +0.14%  '
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.543'
$ws.Range("E45").Value = '  +0.96%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4956'
$ws.Range("E46").Value = '  +0.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.53'
$ws.Range("E47").Value = '  +1.42%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.012'
$ws.Range("E48").Value = '  +1.16%  '
$ws.Range("E49").Value = '  -0.13%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '103.89'
$ws.Range("E50").Value = '  +0.90%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06328'
$ws.Range("E51").Value = '  +0.69%  '
